# Updated legacy GSC export data:
# The oldest date row (2025-10-29) in the "Chart" sheet's coverage table
# has aged out of the rolling export window and must be removed. Deleting
# the row shifts every subsequent date row up by one, which is exactly
# what the refreshed export reflects (last row 2026-01-24 drops off too).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Chart" sheet holds the Date/Not indexed/Indexed/Impressions table

# Row 2 is the first data row (A2 = "2025-10-29"); delete it and shift the
# remaining rows up, matching the new top-of-range date "2025-10-30".
$ws.Rows.Item(2).Delete()
